$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Header cell F1: "disponible" - bold, centered horizontally, top vertical alignment
$f1 = $ws.Range("F1")
$f1.Value = "disponible"
$f1.NumberFormat = "@"
$f1.Font.Bold = $true
$f1.HorizontalAlignment = -4108
$f1.VerticalAlignment = -4160

# Body cells F2:F15: "1" as text, centered horizontally, top vertical alignment
$body = $ws.Range("F2:F15")
$body.NumberFormat = "@"
$body.HorizontalAlignment = -4108
$body.VerticalAlignment = -4160
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 6).Value = "1"
}

# Update the active selection to F1
$f1.Select() | Out-Null
